$wb = $excel.ActiveWorkbook

# Map of sheet name -> (cell, old id) so we can replace id='...' with class='...'
$targets = @(
    @{ Sheet = "!!Main root";          Cell = "A2"; Id = "MainRoot" },
    @{ Sheet = "!!Nodes";              Cell = "A1"; Id = "Node" },
    @{ Sheet = "!!Node friends";       Cell = "A1"; Id = "NodeFriend" },
    @{ Sheet = "!!Leaves";             Cell = "A1"; Id = "Leaf" },
    @{ Sheet = "!!One to many rows";   Cell = "A1"; Id = "OneToManyRow" }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    $cell = $ws.Range($t.Cell)
    $oldValue = "!!ObjTables type='Data' id='" + $t.Id + "'"
    $newValue = "!!ObjTables type='Data' class='" + $t.Id + "'"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
